$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S5").Value = 'etw. aus der Hand schütteln; zaubern; etw. aus dem Handgelenk schütteln; zaubern; hat; aus dem Ärmel geschüttelt; aus dem Ärmel schütteln; etw. aus der Hand schütteln / zaubern, etw. aus dem Handgelenk schütteln /;  zaubern'
$ws.Range("S9").Value = 'ein Auge; zudrücken; drücke, zu; zudrückt; beide Augen zudrücken'
$ws.Range("S12").Value = 'habe; Bammel; hat; hatte; Schiss; Schiss haben '
$ws.Range("S13").Value = 'auf die Barrikaden; gehen; gingen; Barrikade '
$ws.Range("S15").Value = 'blau; machen; mache; hat; gemacht; machst; haben; blauen Montag machen; Montag;'
$ws.Range("S16").Value = 'blau; sein; warst; war; blau vor den Augen werden'
$ws.Range("S18").Value = 'am Boden; zerstört; sein; war; ist; waren; den Boden unter den Füßen verlieren'
$ws.Range("S19").Value = 'im; gleichen; einem; Boot; sitzen; to be in the same boat;'
$ws.Range("S20").Value = 'um den heißen Brei; reden; redest; herum; wie die Katze um den heißen Brei schleichen'
$ws.Range("S25").Value = 'guter Dinge; sein; ist; bin; bleibt; aller guten Dinge sind; gute Dinge sein'
$ws.Range("S38").Value = 'alles im grünen Bereich; sein; ist; Alles in einem grünen Bereich; Alles ist grüner Bereich '
$ws.Range("S39").Value = 'einen Kloß im Hals; haben; hatte; bekomme; Knödel; Kloß; Knödeltenor '
$ws.Range("S41").Value = 'aus der Haut; fahren; bin; gefahren; fährt; Nicht aus seiner Haut können'
$ws.Range("S54").Value = 'auf dem Holzweg; sein; war; bist; ist; "Holzweg";'
$ws.Range("S61").Value = 'kalt; machen; macht; “kalt”;'
$ws.Range("S62").Value = 'über einen Kamm; scheren; Bader;'
$ws.Range("S64").Value = 'einen Kater; haben; hatte; Katers; bin; verkatert; sein'
$ws.Range("S65").Value = 'die Katze aus dem Sack; lassen; ist; los; Katze, lässt; die Katze im Sack lassen'
$ws.Range("S68").Value = 'die Klappe; halten; halt; Klappe;'
$ws.Range("S73").Value = 'einen Korb; geben; haben; gegeben; gibt; hat; bekommen;'
$ws.Range("S76").Value = 'bei; Laune; halten; hält; Luna'
$ws.Range("S78").Value = 'den Löffel; abgeben; abgegeben; hat; Löffel;'
$ws.Range("S83").Value = 'einen Narren; an; gefressen; haben; hast; hat; ist vernarrt in; vernarrt in jmdn. / etw. sein;'
$ws.Range("S86").Value = 'Nerven; liegen; blank; lagen; die Nerven; behalten;'
$ws.Range("S89").Value = 'auf die Palme; bringen; bringt; komm; von der Palme; herunter; von der Palme wieder herunterkommen;'
$ws.Range("S97").Value = 'ein totes Pferd; reiten; wird; hat; geritten; riding a dead horse'
$ws.Range("S105").Value = 'im; in; eigenen Saft; schmoren; wirst; schmore; schmort; jmdn. im eigenen Saft schmoren lassen'
$ws.Range("S108").Value = 'schwarzes; das schwarze; Schaf; sein; ist; war; Sündenbock sein, Buhmann sein'
$ws.Range("S111").Value = 'eine; Schlaftablette; sein; ist; Schlaftabletten; sind; dröge sein;'
$ws.Range("S114").Value = 'Scnhee von gestern; olle Kamelle; kalter Kaffee'
$ws.Range("S117").Value = 'kein; armes; Schwein'
$ws.Range("S118").Value = 'Schwein; haben; gehabt; hatten; Sau; Schützenfeste;'
$ws.Range("S121").Value = 'an die; Substanz; gehen; ging; geht; '
$ws.Range("S124").Value = 'nicht alle Tassen im Schrank; haben; hat; habt; eine Macke haben; einen Dachschaden haben; eine Schraube locker haben; nicht alle Latten am Zaun haben; nicht ganz dicht sein;'
$ws.Range("S129").Value = 'einen Vogel; haben, hat; hast; zeigen; zeigte; jmdm. einen Vogel zeigen'
$ws.Range("S130").Value = 'einen an der; Waffel; haben; hast; waffeln; waffle'
$ws.Range("S132").Value = 'ins kalte Wasser; springen; jmdn. ins kalte Wasser werfen/schmeißen;'
$ws.Range("S133").Value = 'nah am Wasser; gebaut; sein; bin; war; ist; Heulsuse; '
$ws.Range("S139").Value = 'sich; fühlen; wie durch den Wolf; gedreht; fühle; mich; durch den; Fleischwolf'
$ws.Range("S142").Value = 'sich; für; ins Zeug; legen; dich; uns; sich für jmdn. ins Zeug legen'
